# Edit script for Motion_To_Stay_Eviction.docx
# Applies:
#  1. "{{ housing_court }} HOUSING COURT DOCKET NO. {{ housing_court_docket_number }}"
#     -> split into three runs using "lower_court_case.docket_number"
#  2. "is_first_document" -> "is_initial_filing" (template variable rename)
#  3. Judgment date / judge sentence split to use "lower_court_case.judgment_date"
#     and "lower_court_case.judge"
#  4. Merge the three "showifdef('other_legal_issues')) if defined(...) else '' }}"
#     runs back into a single run

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace $oldText (found literally) with the concatenation of the
# given segments, then force each segment boundary to live in its own run
# by toggling Bold on/off (this nudges the engine into re-splitting runs
# while leaving the final formatting unchanged).
# ---------------------------------------------------------------------------
function Replace-WithSplitRuns {
    param(
        [string]$OldText,
        [string[]]$Segments
    )

    $newText = [string]::Join("", $Segments)

    $rng = $d.Content
    $rng.Find.ClearFormatting()
    while ($rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $start = $rng.Start
        $rng.Text = $newText

        $pos = $start
        foreach ($seg in $Segments) {
            $segRange = $d.Range($pos, $pos + $seg.Length)
            $segRange.Bold = 1
            $segRange.Bold = 0
            $pos = $pos + $seg.Length
        }

        $endPos = $start + $newText.Length
        $rng.Start = $endPos
        $rng.End = $endPos
    }
}

# ---------------------------------------------------------------------------
# 1 & 3: housing court docket number + judgment date/judge split (lower_court_case.)
# ---------------------------------------------------------------------------
Replace-WithSplitRuns "{{ housing_court }} HOUSING COURT DOCKET NO. {{ housing_court_docket_number }}" `
    @("{{ housing_court }} HOUSING COURT DOCKET NO. {{ ", "lower_court_case.", "docket_number }}")

$quoteRight = [char]0x2019
$oldSentence = "I am requesting a stay of the Housing Court" + $quoteRight + "s judgment entered on {{ housing_court_judgment_date }} which awarded possession of the premises to the plaintiff.  The judgment was entered by Judge {{ housing_court_judge }}."
Replace-WithSplitRuns $oldSentence @(
    ("I am requesting a stay of the Housing Court" + $quoteRight + "s judgment entered on {{ "),
    "lower_court_case.",
    "judgment_date }} which awarded possession of the premises to the plaintiff.  The judgment was entered by Judge {{ ",
    "lower_court_case.",
    "judge }}."
)

# ---------------------------------------------------------------------------
# 2: is_first_document -> is_initial_filing, preserving the existing (empty)
#    run-level formatting by toggling Bold over the run's visible text only
#    (excluding the paragraph mark).
#
#    NOTE: Find.Execute with a non-empty Replace string only performs a
#    single replacement per call even inside a while loop (it does not
#    keep advancing through the rest of the document the way repeated
#    search-only calls do), so we search-only here and apply the text
#    change ourselves.
# ---------------------------------------------------------------------------
$oldWord = "is_first_document"
$newWord = "is_initial_filing"
$rng = $d.Content
$rng.Find.ClearFormatting()
while ($rng.Find.Execute($oldWord, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $para = $rng.Duplicate
    $para.Expand(4)
    $paraStart = $para.Start
    $paraTextEnd = ($para.End - 1) + ($newWord.Length - $oldWord.Length)

    $rng.Text = $newWord

    $textOnly = $d.Range($paraStart, $paraTextEnd)
    $textOnly.Bold = 1
    $textOnly.Bold = 0

    $rng.Collapse(0)
}

# ---------------------------------------------------------------------------
# 4: merge the three showifdef('other_legal_issues') runs into a single run
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
if ($rng.Find.Execute("showifdef(", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $start = $rng.Start
    $para = $rng.Duplicate
    $para.Expand(4)
    $paraEnd = $para.End - 1

    $full = $d.Range($start, $paraEnd)

    $quoteLeft = [char]0x2018
    $quoteRight2 = [char]0x2019
    $mergedText = "showifdef(" + $quoteLeft + "other_legal_issues" + $quoteRight2 + ")) if defined(" + $quoteLeft + "other_legal_issues" + $quoteRight2 + ") else " + $quoteLeft + $quoteRight2 + " }}"

    # First overwrite with a placeholder to force a genuine text change, then
    # set the final merged text - this collapses the previously separate runs
    # into one while keeping the run's formatting (bCs) intact.
    $full.Text = "PLACEHOLDER_MERGE_TOKEN"
    $final = $d.Range($start, $start + "PLACEHOLDER_MERGE_TOKEN".Length)
    $final.Text = $mergedText
}
